$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.888.54'
$ws.Range("E2").Value = '  +2.99%  '
$ws.Range("D3").Value = '3.560.07'
$ws.Range("E3").Value = '  +1.98%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.75'
$ws.Range("E5").Value = '  +5.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.24'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  +2.18%  '
$ws.Range("D8").Value = '3.555.49'
$ws.Range("E8").Value = '  +1.99%  '
$ws.Range("E9").Value = '  -0.13%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.198'
$ws.Range("E10").Value = '  +5.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.98'
$ws.Range("E11").Value = '  +6.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.587'
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.78'
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000278'
$ws.Range("E14").Value = '  +2.79%  '
$ws.Range("D15").Value = '4.131.07'
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.40'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '619.27'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '70.923.94'
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("D19").Value = '3.543.69'
$ws.Range("E19").Value = '  +1.07%  '
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.38'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.883'
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.61'
$ws.Range("E23").Value = '  -12.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.73'
$ws.Range("E24").Value = '  -0.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '96.87'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.62'
$ws.Range("E28").Value = '  +0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.44'
$ws.Range("E29").Value = '  +3.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.09'
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.52'
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("E32").Value = '  -1.88%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.03'
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.31'
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '574.27'
$ws.Range("E35").Value = '  -8.35%  '
$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.64'
$ws.Range("E36").Value = '  +4.95%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.101'
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.83'
$ws.Range("E38").Value = '  +1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '57.70'
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("E40").Value = '  +6.01%  '
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("E42").Value = '  +4.94%  '
$ws.Range("D43").Value = '3.364.62'
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("E44").Value = '  -1.41%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.01'
$ws.Range("E45").Value = '  +8.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.09'
$ws.Range("E46").Value = '  +1.38%  '
$ws.Range("D47").Value = '0.0₃0708'
$ws.Range("E47").Value = '  +2.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.64'
$ws.Range("E48").Value = '  +3.65%  '
$ws.Range("E49").Value = '  +0.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.94'
$ws.Range("E50").Value = '  +1.96%  '
$ws.Range("E51").Value = '  +1.92%  '
